$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1: "Valeurs reelles" ---
# Fill diagonal values that were previously blank
$ws1.Range("E23").Value = 1.404999971389771
$ws1.Range("D24").Value = 1.404999971389771
$ws1.Range("C25").Value = 1.404999971389771

# New row 26
$ws1.Range("A26").Value = 45684
$ws1.Range("A26").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B26").Value = 1.404999971389771

# --- Sheet 2: "Predictions" ---
$ws2.Range("A26").Value = 45684
$ws2.Range("A26").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B26").Value = 1.402897357940674
$ws2.Range("C26").Value = 1.373286366462708
$ws2.Range("D26").Value = 1.338665127754211
